function Set-TextValue($Ws, $Ref, $Val) {
    $cell = $Ws.Range($Ref)
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $Val
    $cell.Style = $style
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
Set-TextValue $ws "D2" "244.34"
Set-TextValue $ws "D3" "21.86"
Set-TextValue $ws "D4" "5.391"
Set-TextValue $ws "D6" "3.391"
Set-TextValue $ws "D7" "0.8180"
Set-TextValue $ws "D8" "0.9546"
Set-TextValue $ws "D10" "0.07437"
Set-TextValue $ws "D11" "0.03291"
Set-TextValue $ws "D12" "0.03054"
Set-TextValue $ws "D13" "0.09408"
Set-TextValue $ws "D14" "4.007"
Set-TextValue $ws "D15" "0.001590"
Set-TextValue $ws "D16" "0.04796"
Set-TextValue $ws "D18" "0.005421"
Set-TextValue $ws "D19" "0.004153"
Set-TextValue $ws "D20" "0.0009911"
Set-TextValue $ws "D21" "3.672"
Set-TextValue $ws "D22" "6.422"
Set-TextValue $ws "D23" "2.189"
Set-TextValue $ws "D40" "0.03994"

# --- Rows 41-43 reshuffled (KickToken/BKEXToken/CEJI rotate) ---
Set-TextValue $ws "B41" "BKEXToken"
Set-TextValue $ws "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D41" "0.1075"
Set-TextValue $ws "E41" "40BKEXTokenBKK"
Set-TextValue $ws "B42" "CEJI"
Set-TextValue $ws "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.002721"
Set-TextValue $ws "E42" "41CEJICEJI"
Set-TextValue $ws "B43" "KickToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D43" "0.003040"
Set-TextValue $ws "E43" "42KickTokenKICK"

# --- Misc single-cell updates ---
Set-TextValue $ws "D44" "0.005822"
Set-TextValue $ws "E47" "46CoinbaseStockTokenCOINBestin24h"
Set-TextValue $ws "D48" "0.004429"
